$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 values (product attributes changed)
$ws.Range("A2").Value = "red shirts"
$ws.Range("B2").Value = "MSP64LT20025"
$ws.Range("C2").Value = "Navy"
$ws.Range("D2").Value = "M"

# Add a new "Product name" column
$ws.Range("E1").Value = "Product name"
$ws.Range("E2").Value = "AWEARNESS Kenneth Cole"
$ws.Columns.Item(5).ColumnWidth = 22.5

# Match the selection left by the author
$ws.Range("D2").Select()
